$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = 46009

$ws.Range("B37").Value = 67

$ws.Range("A37:B37").Select()
